$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(
    @("2021-06-19", "overview", "K02000001", "United Kingdom", 4620968, 10321, 14, 127970),
    @("2021-06-20", "overview", "K02000001", "United Kingdom", 4630040, 9284, 6, 127976),
    @("2021-06-21", "overview", "K02000001", "United Kingdom", 4640507, 10633, 5, 127981),
    @("2021-06-22", "overview", "K02000001", "United Kingdom", 4651988, 11625, 27, 128008),
    @("2021-06-23", "overview", "K02000001", "United Kingdom", 4667870, 16135, 19, 128027),
    @("2021-06-24", "overview", "K02000001", "United Kingdom", 4684572, 16703, 21, 128048),
    @("2021-06-25", "overview", "K02000001", "United Kingdom", 4699868, 15810, 18, 128066)
)

$startRow = 312
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Column A holds a date-like string ("YYYY-MM-DD") that must stay plain
    # text (as in the source data) rather than be auto-converted to a date
    # serial by Excel's smart-entry parsing. Force text formatting for the
    # entry, then drop the style back to Normal so the cell keeps the sheet's
    # default (unstyled) formatting, matching the other rows.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $data[0]
    $cellA.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
}
